# Applies a weekly re-shuffle of the "Fecha", "Volumen" and price columns
# (D, J, K, L, M, P) across several rows of the sheet, as described by the
# commit "Fruta / hortaliza, semanal".
#
# The rows whose data moves are linked by the permutation:
#   2 <- 16, 16 <- 6, 6 <- 13, 13 <- 10, 10 <- 2      (cycle 1)
#   3 <- 11, 11 <- 4, 4 <- 15, 15 <- 17, 17 <- 12, 12 <- 8, 8 <- 3   (cycle 2)
# i.e. row R receives the values that used to belong to row Source(R).
# Rows 5, 7, 9, 14 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a group for each row.
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Mapping: destination row -> source row (values read from source, written to destination)
$map = @{
    2  = 16
    16 = 6
    6  = 13
    13 = 10
    10 = 2
    3  = 11
    11 = 4
    4  = 15
    15 = 17
    17 = 12
    12 = 8
    8  = 3
}

# Snapshot the original values for every row involved, BEFORE writing anything,
# since several rows both give and receive data (cycles).
# NOTE: use Value2 (not Value) for reads -- in this runtime, .Value alone
# does not reliably return the scalar, while .Value2 does.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the shuffled values using the snapshot (so later writes don't
# clobber data still needed for other rows).
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
